$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 12; $row++) {
    $ws.Cells.Item($row, 2).Value = "2020-05-09 16:47"
    $ws.Cells.Item($row, 3).Value = 23
}
